$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the computed percentages (row 2 / row 3, column C) with the
# corrected figures from the regrouped match_on_abstract results.
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "1.93"
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "99.17"

# Re-style the whole table with the report's body font (Atkinson
# Hyperlegible Regular, 10pt) and drop the old thin-box border grid.
$all = $ws.Range("A1:C3")
$all.Font.Name = "Atkinson Hyperlegible Regular"
$all.Font.Size = 10
$all.Borders.LineStyle = -4142

# Select the whole table, matching the resaved worksheet view.
$ws.Range("A1:C3").Select() | Out-Null
